$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date value in A1 (keeps existing date format / style)
$ws.Range("A1").Value = 45436

# Update the price list values in column D
$ws.Range("D32").Value = 219.124
$ws.Range("D33").Value = 313.543
$ws.Range("D34").Value = 417.457
$ws.Range("D35").Value = 429.794
$ws.Range("D36").Value = 563.266
$ws.Range("D37").Value = 644.069
$ws.Range("D38").Value = 771.267
$ws.Range("D39").Value = 918.41
